$d = $word.ActiveDocument

# Items 6, 7 and 8 keep the same visible text, but the diff shows their
# runs being coalesced into fewer runs (an artifact of Word re-saving the
# paragraphs after the nearby edit). Re-asserting each paragraph's own
# text via Find & Replace reproduces that same-format run coalescing.

$item6tail = "squared). Show a dialog box or your code here."
$d.Content.Find.Execute($item6tail, $true, $false, $false, $false, $false, $true, 1, $false, $item6tail, 2)

$item7 = "7. Recode the sex data into a different variable with 0 = Male, 1 = Female. Show a dialog box or your code here."
$d.Content.Find.Execute($item7, $true, $false, $false, $false, $false, $true, 1, $false, $item7, 2)

$item8tail = " and all the females appear second in the data. Show a dialog box or your code here. Show a listing of the sorted data."
$d.Content.Find.Execute($item8tail, $true, $false, $false, $false, $false, $true, 1, $false, $item8tail, 2)

# Item 9: replace everything after "9. " with the new remark.
$old9 = "Create a data dictionary for this file. The data dictionary should specify who created the file (Monica Gaddis), when it was created (date unknown), how many rows and columns are in the dataset, what format the data was originally stored in (text file with comma delimiters), the variable names for each column of data (except bmi) with a brief description of the variable. Be sure to include units of measurement and the categories associated with any number codes. Copy the information from your data dictionary here."
$new9 = "Ignore this question. I am moving it to the homework on the next module."
$d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
